$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.896.79"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.879.97"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "333.19"
$ws.Range("E5").Value = "  +3.84%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.4730"
$ws.Range("E7").Value = "  +5.96%  "
$ws.Range("D8").Value = "0.3977"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("D9").Value = "48.73"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "0.08062"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").Value = "1.028"
$ws.Range("E11").Value = "  +1.86%  "
$ws.Range("D12").Value = "21.95"
$ws.Range("E12").Value = "  +2.99%  "
$ws.Range("D13").Value = "1.915.82"
$ws.Range("E13").Value = "  +4.34%  "
$ws.Range("D14").Value = "5.961"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "7.192"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "87.26"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "0.06626"
$ws.Range("E19").Value = "  +1.88%  "
$ws.Range("D20").Value = "17.29"
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").Value = "28.050.91"
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").Value = "  +2.84%  "
$ws.Range("D25").Value = "2.306"
$ws.Range("E25").Value = "  +2.59%  "
$ws.Range("D26").Value = "2.118.08"
$ws.Range("E26").Value = "  +2.99%  "
$ws.Range("D27").Value = "157.19"
$ws.Range("E27").Value = "  +3.84%  "
$ws.Range("D28").Value = "20.23"
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("D29").Value = "2.106"
$ws.Range("E29").Value = "  +2.74%  "
$ws.Range("D30").Value = "5.629"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").Value = "0.9802"
$ws.Range("E32").Value = "  +5.94%  "
$ws.Range("D33").Value = "0.09565"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").Value = "1.463"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").Value = "5.327"
$ws.Range("E36").Value = "  +1.90%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.02264"
$ws.Range("E37").Value = "  +2.61%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06109"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("D39").Value = "1.228"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").Value = "8.227"
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("D41").Value = "0.6035"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "0.1904"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("D44").Value = "10.32"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").Value = "0.5740"
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("D48").Value = "3.415"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").Value = "1.947"
$ws.Range("D50").Value = "0.06829"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "114.07"
$ws.Range("E51").Value = "  +5.49%  "
